# Adds a new "2022-Q4" sheet (with its own fund-holdings table) right after
# the "总计" (Total) sheet, inserts a matching summary row into "总计", and
# leaves the pre-existing "2021-Q3" sheet's data intact (it just shifts
# from position 2 to position 3, unchanged).
#
# To keep the engine's internal sheetId numbering lined up with the target
# (总计=1, 2022-Q4=2, 2021-Q3=3) we duplicate the *existing* "2021-Q3" sheet
# in place (the duplicate inherits the next free sheetId and keeps the
# "2021-Q3" name/data going forward), then repurpose the original sheet
# object - still holding sheetId 2 - as the new "2022-Q4" sheet.

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item(1)      # "总计"
$oldQ  = $wb.Worksheets.Item("2021-Q3")

# ---------------------------------------------------------------------
# 1) Duplicate "2021-Q3" right after itself; the duplicate carries the
#    "2021-Q3" name/data forward. The original sheet object gets reused
#    (and renamed) for the new "2022-Q4" content.
# ---------------------------------------------------------------------
$oldQ.Copy($null, $oldQ)
$dup = $wb.Worksheets.Item(3)

$oldQ.Name = "2022-Q4"
$dup.Name  = "2021-Q3"

$q4 = $oldQ

# ---------------------------------------------------------------------
# 2) Wipe the old fund table out of "2022-Q4" and lay down the new one.
# ---------------------------------------------------------------------
$q4.Cells.Clear()

$total.Range("B1:D1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats

$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

$total.Range("A2").Copy()
$q4.Range("A2").PasteSpecial(-4122)      # xlPasteFormats
$q4.Range("A2").Value = 0

# B2 is a numeric-looking fund code ("004685") that must keep its leading
# zero, i.e. stay text, not become the number 4685. D2,E2,F2,G2 are the
# same situation (numeric-looking text in the source data). Force text via
# the quote-prefix so they keep their original string form.
$q4.Range("B2").NumberFormat = "@"
$q4.Range("B2").Value = "004685"
$q4.Range("C2").Value = "金元顺安元启灵活配置混合"

$q4.Range("D2:G2").NumberFormat = "@"
$q4.Range("D2").Value = "15.29"
$q4.Range("E2").Value = "76.11"
$q4.Range("F2").Value = "1.00"
$q4.Range("G2").Value = "0.1529"

$q4.Range("H2").Value = 3

# ---------------------------------------------------------------------
# 3) Push "总计"'s existing 2021-Q3 row down to row 3, then write the new
#    2022-Q4 summary row into row 2.
# ---------------------------------------------------------------------
$total.Range("A2").Copy()
$total.Range("A3").PasteSpecial(-4122)   # xlPasteFormats
$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q3"
$total.Range("C3").Value = 2
$total.Range("D3").Value = 0.14

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.15

# ---------------------------------------------------------------------
# 4) Keep "2021-Q3" as the selected/active tab (it is otherwise untouched
#    by this change).
# ---------------------------------------------------------------------
$dup.Activate()
